$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.437.57'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.514.16'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '540.74'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.60'
$ws.Range('E6').Value = '  -3.80%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  -1.57%  '
$ws.Range('D9').Value = '2.516.07'
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('E12').Value = '  -3.60%  '
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = '2.962.56'
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.39'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = '59.334.94'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '2.513.82'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.10'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '325.20'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.84'
$ws.Range('E23').Value = '  -1.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.34'
$ws.Range('E24').Value = '  +1.62%  '
$ws.Range('E26').Value = '  +2.19%  '
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.83'
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('E29').Value = '  +1.95%  '
$ws.Range('D30').Value = '0.0₃0780'
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.80'
$ws.Range('E31').Value = '  -1.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.03'
$ws.Range('E32').Value = '  +3.07%  '
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.12'
$ws.Range('E34').Value = '  -6.73%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.44'
$ws.Range('E35').Value = '  -2.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.50'
$ws.Range('E36').Value = '  -1.33%  '
$ws.Range('E37').Value = '  -2.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.60'
$ws.Range('E38').Value = '  -0.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.91'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('E40').Value = '  -0.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.811'
$ws.Range('E41').Value = '  -2.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.25'
$ws.Range('E42').Value = '  -5.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '281.94'
$ws.Range('E43').Value = '  -4.60%  '
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.88'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.598'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.71'
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('E50').Value = '  -1.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.90'
$ws.Range('E51').Value = '  -1.98%  '
